$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.631.78"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.18"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.38"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.94"
$ws.Range("E8").Value = "  +5.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0587"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.788.02"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.561.77"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.649.35"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.46"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.62"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0682"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.81"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.106"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.80"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.25"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.405.72"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.517"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0464"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.768"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.91"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.699.76"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.866"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.82"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.72"
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").Value = "  -0.43%  "
